$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices: A=1, B=2, C=3, D=4, E=5

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).Value = "71.143.76"
$ws.Cells.Item(2, 5).Value = "  +2.87%  "

# Row 3 - Ethereum
$ws.Cells.Item(3, 4).Value = "3.806.46"
$ws.Cells.Item(3, 5).Value = "  +0.93%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 5).Value = "  +0.05%  "

# Row 5 - BNB
$ws.Cells.Item(5, 4).Value = "698.13"
$ws.Cells.Item(5, 5).Value = "  +10.61%  "

# Row 6 - Solana
$ws.Cells.Item(6, 4).Value = "172.92"
$ws.Cells.Item(6, 5).Value = "  +3.70%  "

# Row 7 - LidoStakedEther
$ws.Cells.Item(7, 4).Value = "3.805.08"
$ws.Cells.Item(7, 5).Value = "  +0.95%  "

# Row 8 - USDC
$ws.Cells.Item(8, 5).Value = "  -0.06%  "

# Row 9 - XRP
$ws.Cells.Item(9, 4).Value = "0.525"
$ws.Cells.Item(9, 5).Value = "  +0.82%  "

# Row 10 - Dogecoin
$ws.Cells.Item(10, 5).Value = "  +2.80%  "

# Row 11 - Toncoin
$ws.Cells.Item(11, 4).Value = "7.56"
$ws.Cells.Item(11, 5).Value = "  +12.10%  "

# Row 12 - Cardano
$ws.Cells.Item(12, 5).Value = "  +0.55%  "

# Row 13 - ShibaInu
$ws.Cells.Item(13, 4).Value = "0.0000257"
$ws.Cells.Item(13, 5).Value = "  +7.82%  "

# Row 14 - Avalanche
$ws.Cells.Item(14, 4).Value = "36.22"
$ws.Cells.Item(14, 5).Value = "  +3.39%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Cells.Item(15, 4).Value = "4.448.69"
$ws.Cells.Item(15, 5).Value = "  +0.94%  "

# Row 16 - WrappedEther
$ws.Cells.Item(16, 4).Value = "3.812.15"
$ws.Cells.Item(16, 5).Value = "  +0.83%  "

# Row 17 - WrappedBTC
$ws.Cells.Item(17, 4).Value = "71.146.87"
$ws.Cells.Item(17, 5).Value = "  +2.83%  "

# Row 18 - Chainlink
$ws.Cells.Item(18, 4).Value = "17.85"
$ws.Cells.Item(18, 5).Value = "  +1.43%  "

# Row 19 - Polkadot
$ws.Cells.Item(19, 5).Value = "  +3.10%  "

# Row 20 - TRON
$ws.Cells.Item(20, 5).Value = "  +1.20%  "

# Row 21 - Uniswap
$ws.Cells.Item(21, 5).Value = "  +17.66%  "

# Row 22 - BitcoinCash
$ws.Cells.Item(22, 4).Value = "482.16"
$ws.Cells.Item(22, 5).Value = "  +4.17%  "

# Row 23 - Polygon
$ws.Cells.Item(23, 4).Value = "0.716"
$ws.Cells.Item(23, 5).Value = "  +1.47%  "

# Row 24 - Litecoin
$ws.Cells.Item(24, 4).Value = "83.85"
$ws.Cells.Item(24, 5).Value = "  +2.18%  "

# Row 25 - PEPE
$ws.Cells.Item(25, 4).Value = "0.0000146"
$ws.Cells.Item(25, 5).Value = "  +0.87%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Cells.Item(26, 4).Value = "12.37"
$ws.Cells.Item(26, 5).Value = "  +2.42%  "

# Row 27 - was Fetch.AI, now RenderToken
$ws.Cells.Item(27, 2).Value = "RenderToken"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(27, 4).Value = "10.51"
$ws.Cells.Item(27, 5).Value = "  +4.59%  "

# Row 28 - was RenderToken, now Fetch.AI
$ws.Cells.Item(28, 2).Value = "Fetch.AI"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(28, 4).Value = "2.17"
$ws.Cells.Item(28, 5).Value = "  +1.37%  "

# Row 29 - WrappedeETH
$ws.Cells.Item(29, 4).Value = "3.957.31"
$ws.Cells.Item(29, 5).Value = "  +0.92%  "

# Row 30 - Dai
$ws.Cells.Item(30, 5).Value = "  +0.01%  "

# Row 31 - PancakeSwap
$ws.Cells.Item(31, 4).Value = "3.07"
$ws.Cells.Item(31, 5).Value = "  +14.23%  "

# Row 32 - ImmutableX
$ws.Cells.Item(32, 5).Value = "  -0.31%  "

# Row 33 - NEARProtocol
$ws.Cells.Item(33, 4).Value = "7.57"
$ws.Cells.Item(33, 5).Value = "  +6.98%  "

# Row 34 - EthereumClassic
$ws.Cells.Item(34, 4).Value = "29.59"
$ws.Cells.Item(34, 5).Value = "  +3.86%  "

# Row 35 - Kaspa
$ws.Cells.Item(35, 4).Value = "0.177"
$ws.Cells.Item(35, 5).Value = "  -1.06%  "

# Row 36 - Aptos
$ws.Cells.Item(36, 4).Value = "9.21"
$ws.Cells.Item(36, 5).Value = "  +3.04%  "

# Row 37 - Binance-PegBSC-USD
$ws.Cells.Item(37, 4).Value = "1.00"
$ws.Cells.Item(37, 5).Value = "  +0.07%  "

# Row 38 - RenzoRestakedETH
$ws.Cells.Item(38, 4).Value = "3.756.08"
$ws.Cells.Item(38, 5).Value = "  +0.85%  "

# Row 39 - Hedera
$ws.Cells.Item(39, 5).Value = "  +1.54%  "

# Row 40 - dogwifhat
$ws.Cells.Item(40, 5).Value = "  +5.69%  "

# Row 41 - Filecoin
$ws.Cells.Item(41, 5).Value = "  +3.12%  "

# Row 42 - Stacks
$ws.Cells.Item(42, 5).Value = "  +11.30%  "

# Row 43 - was Mantle, now FLOKI
$ws.Cells.Item(43, 2).Value = "FLOKI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Cells.Item(43, 4).Value = "0.000327"
$ws.Cells.Item(43, 5).Value = "  +22.69%  "

# Row 44 - was FirstDigitalUSD, now Mantle
$ws.Cells.Item(44, 2).Value = "Mantle"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(44, 4).Value = "0.970"
$ws.Cells.Item(44, 5).Value = "  +0.66%  "

# Row 45 - was FLOKI, now FirstDigitalUSD
$ws.Cells.Item(45, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(45, 4).Value = "1.00"
$ws.Cells.Item(45, 5).Value = "  +0.07%  "

# Row 47 - Arweave
$ws.Cells.Item(47, 4).Value = "45.58"
$ws.Cells.Item(47, 5).Value = "  +5.05%  "

# Row 48 - OKB
$ws.Cells.Item(48, 4).Value = "49.33"

# Row 49 - Monero
$ws.Cells.Item(49, 4).Value = "160.19"
$ws.Cells.Item(49, 5).Value = "  +1.43%  "

# Row 50 - ONDO
$ws.Cells.Item(50, 5).Value = "  -1.81%  "

# Row 51 - TheGraph
$ws.Cells.Item(51, 5).Value = "  +1.33%  "
